# Apply the "repull data" edit: update the dSF ("F") column values for a
# specific set of rows on the active sheet to reflect freshly re-pulled data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> new value for column F (dSF)
$updates = @{
    4  = -5
    12 = -3
    13 = 1
    23 = -3
    28 = 1
    35 = -8
    41 = -3
    42 = 6
    44 = -4
    45 = -1
    46 = -2
    54 = 1
    63 = -4
    65 = 2
    67 = 1
    74 = -1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
